$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.785.78"
$ws.Range("D3").Value = "2.101.09"
$ws.Range("E3").Value = "  +0.43%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.09"
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("E6").Value = "  +0.98%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.65"
$ws.Range("E7").Value = "  +1.61%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  +2.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0845"
$ws.Range("E10").Value = "  +0.37%  "
$ws.Range("E11").Value = "  +0.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.47"
$ws.Range("E12").Value = "  +5.71%  "
$ws.Range("D13").Value = "2.412.43"
$ws.Range("E13").Value = "  +0.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.07"
$ws.Range("E14").Value = "  -0.32%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "2.312.74"
$ws.Range("E15").Value = "  +9.46%  "
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.809"
$ws.Range("E16").Value = "  +4.53%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.50"
$ws.Range("E17").Value = "  +0.56%  "
$ws.Range("D18").Value = "38.801.38"
$ws.Range("E18").Value = "  +1.95%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.00"
$ws.Range("E19").Value = "  +2.70%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.09"
$ws.Range("E20").Value = "  +1.39%  "
$ws.Range("D21").Value = "0.0₃0841"
$ws.Range("E21").Value = "  +0.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.64"
$ws.Range("E22").Value = "  +1.67%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("E24").Value = "  -2.28%  "
$ws.Range("E25").Value = "  +0.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "171.71"
$ws.Range("E26").Value = "  +1.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.54"
$ws.Range("E27").Value = "  +1.19%  "
$ws.Range("E28").Value = "  +5.95%  "
$ws.Range("E29").Value = "  +4.97%  "
$ws.Range("E30").Value = "  +2.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.48"
$ws.Range("E31").Value = "  +3.87%  "
$ws.Range("E32").Value = "  +1.22%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.54"
$ws.Range("E33").Value = "  +2.54%  "
$ws.Range("E34").Value = "  +1.68%  "
$ws.Range("E35").Value = "  +2.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.52"
$ws.Range("E36").Value = "  +1.89%  "
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("E38").Value = "  +1.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.17"
$ws.Range("E40").Value = "  +0.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0229"
$ws.Range("E41").Value = "  +4.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "101.50"
$ws.Range("E42").Value = "  +1.34%  "
$ws.Range("D43").Value = "1.532.77"
$ws.Range("E43").Value = "  -1.55%  "
$ws.Range("E44").Value = "  -0.91%  "
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0911"
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.74"
$ws.Range("E46").Value = "  +3.34%  "
$ws.Range("E47").Value = "  +2.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.11"
$ws.Range("E48").Value = "  -0.76%  "
$ws.Range("E49").Value = "  +2.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.97"
$ws.Range("E50").Value = "  -0.50%  "
$ws.Range("D51").Value = "2.295.19"
$ws.Range("E51").Value = "  +0.31%  "
